$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the "Categoria" column (2nd column) entirely.
$t.Columns.Item(2).Delete()

# Update the "n" row value.
$t.Cell(2, 2).Range.Text = "426"

# For each medication, the former "NAO" row becomes the surviving row
# (its label is updated to "<Drug> = SIM (%)" and its value cell takes
# the percentage that used to live in the following "SIM" row), and the
# old "SIM" row is deleted.
$t.Cell(3, 1).Range.Text = "Dabigatrana = SIM (%)"
$t.Cell(3, 2).Range.Text = "5 ( 1.2)"
$t.Rows.Item(4).Delete()

$t.Cell(4, 1).Range.Text = "Enoxaparina = SIM (%)"
$t.Cell(4, 2).Range.Text = "384 (93.2)"
$t.Rows.Item(5).Delete()

$t.Cell(5, 1).Range.Text = "Rivoraxabana = SIM (%)"
$t.Cell(5, 2).Range.Text = "87 (21.0)"
$t.Rows.Item(6).Delete()

$t.Cell(6, 1).Range.Text = "Warfarina = SIM (%)"
$t.Cell(6, 2).Range.Text = "196 (52.0)"
$t.Rows.Item(7).Delete()
